$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 289
$ws.Range('A289').Value = '2025-04-18 10:02'
$ws.Range('B289').Value = 'http://www.scpc.gov.cn/group3/M00/08/FD/rBUtImLCOE-ACmH0AADF69H66dc55.docx'
$ws.Range('C289').Value = 'http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=dc06e261b7f14ee0bf15f2175a30cb24&type=0'
$ws.Range('D289').Value = '抓住关键少数'
$ws.Range('E289').Value = '抓住“关键少数”'
$ws.Range('F289').Value = '平昌县纪委2018年部门决算公开编制说明.docx'
$ws.Range('G289').Value = 'http://www.scpc.gov.cn/public/6601841/12823521.html'
$ws.Range('H289').Value = '平昌县纪律检查委员会2018年度部门决算公开'

# Row 290
$ws.Range('A290').Value = '2025-04-18 15:26'
$ws.Range('B290').Value = 'http://www.scpc.gov.cn/group3/M00/08/11/rBUtImIxibiANVSEAAE1Rnt3DRg20.docx'
$ws.Range('C290').Value = 'http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=0171c3b214e94f7d8118f44eb93c9d03&type=0'
$ws.Range('D290').Value = '社会和保障就业'
$ws.Range('E290').Value = '社会保障和就业'
$ws.Range('F290').Value = '平昌县林业局2019年度部门决算公开编制说明.docx.docx'
$ws.Range('G290').Value = 'http://www.scpc.gov.cn/public/6602341/12988791.html'
$ws.Range('H290').Value = '平昌县林业局2019年度部门决算公开'

# Row 291
$ws.Range('A291').Value = '2025-04-22 09:13'
$ws.Range('B291').Value = 'http://www.scpc.gov.cn/group3/M00/0E/5F/rBUtImS14xqAEZeOAABAYgmoXq020.xlsx'
$ws.Range('C291').Value = 'http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=7ba8e15bffb24dffb06ea8ac7a8c0cb2&type=0'
$ws.Range('D291').Value = '中华人民共和国政府信息公开条例》'
$ws.Range('E291').Value = '《中华人民共和国政府信息公开条例》'
$ws.Range('F291').Value = '附件：巴中市平昌县五木镇政府信息主动公开基本目录.xlsx'
$ws.Range('G291').Value = 'http://www.scpc.gov.cn/ztzl/zfxxzdgkjbml/xzbsc/12617841.html'
$ws.Range('H291').Value = '巴中市平昌县五木镇政府信息主动公开基本目录'

# Row 292
$ws.Range('A292').Value = '2025-04-24 15:10'
$ws.Range('B292').Value = 'http://www.scpc.gov.cn/group3/M00/15/64/rBUtImekVsGABHxgACB4yp4VTL4486.pdf'
$ws.Range('C292').Value = 'http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=79733ccd3d67491280ae8f6d0420797a&type=0'
$ws.Range('D292').Value = '粉粹'
$ws.Range('E292').Value = '粉碎'
$ws.Range('F292').Value = '有机肥料技改项目--公示本.pdf'
$ws.Range('G292').Value = 'http://www.scpc.gov.cn/public/6602261/13851381.html'
$ws.Range('H292').Value = '关于2019年6月19日已受理建设项目环评文件公告(有机肥料技改项目)'

# Row 293
$ws.Range('A293').Value = '2025-04-24 15:30'
$ws.Range('B293').Value = 'http://www.scpc.gov.cn/group3/M00/15/B0/rBUtImfY05iAYwg_ABfNmCNZ1no283.pdf'
$ws.Range('C293').Value = 'http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=3063c91c7155498aa44805c71ce2a8c5&type=0'
$ws.Range('D293').Value = '混凝土浇铸'
$ws.Range('E293').Value = '混凝土浇筑'
$ws.Range('F293').Value = '关于2018年10月11日已受理建设项目环评文件公示(平昌县俱兴砖厂迁建项目)'
$ws.Range('G293').Value = 'http://www.scpc.gov.cn/public/6602261/13851441.html'
$ws.Range('H293').Value = '关于2018年10月11日已受理建设项目环评文件公示(平昌县俱兴砖厂迁建项目)'

# Row 294
$ws.Range('A294').Value = '2025-04-28 14:27'
$ws.Range('B294').Value = 'http://www.scpc.gov.cn/group3/M00/15/64/rBUtImekVsGABHxgACB4yp4VTL4486.pdf'
$ws.Range('C294').Value = 'http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=79733ccd3d67491280ae8f6d0420797a&type=0'
$ws.Range('D294').Value = '粉粹'
$ws.Range('E294').Value = '粉碎'
$ws.Range('F294').Value = '有机肥料技改项目--公示本.pdf'
$ws.Range('G294').Value = 'http://www.scpc.gov.cn/public/6602261/13851381.html'
$ws.Range('H294').Value = '关于2019年6月19日已受理建设项目环评文件公告(有机肥料技改项目)'

# Row 295
$ws.Range('A295').Value = '2025-04-28 14:31'
$ws.Range('B295').Value = 'http://www.scpc.gov.cn/group3/M00/12/70/rBUtImb2DMaAGIJqAAGn1E--QX022.xlsx'
$ws.Range('C295').Value = 'http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=d69de281f45242ae8561aeef6f468ecc&type=0'
$ws.Range('D295').Value = '辩认'
$ws.Range('E295').Value = '辨认'
$ws.Range('F295').Value = '附件：第三批巴中市行政处罚五张清单（试行）.xlsx'
$ws.Range('G295').Value = 'http://www.scpc.gov.cn/public/6601961/13964433.html'
$ws.Range('H295').Value = '中共巴中市委全面依法治市委员会办公室关于印发第三批《巴中市行政处罚五张清单（试行）》的通知'

# Row 296
$ws.Range('A296').Value = '2025-04-28 14:47'
$ws.Range('B296').Value = 'http://www.scpc.gov.cn/group3/M00/15/B0/rBUtImfY05iAYwg_ABfNmCNZ1no283.pdf'
$ws.Range('C296').Value = 'http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=3063c91c7155498aa44805c71ce2a8c5&type=0'
$ws.Range('D296').Value = '混凝土浇铸'
$ws.Range('E296').Value = '混凝土浇筑'
$ws.Range('F296').Value = '关于2018年10月11日已受理建设项目环评文件公示(平昌县俱兴砖厂迁建项目)'
$ws.Range('G296').Value = 'http://www.scpc.gov.cn/public/6602261/13851441.html'
$ws.Range('H296').Value = '关于2018年10月11日已受理建设项目环评文件公示(平昌县俱兴砖厂迁建项目)'

# Row 297
$ws.Range('A297').Value = '2025-04-28 14:48'
$ws.Range('B297').Value = 'http://www.scpc.gov.cn/group3/M00/0E/5F/rBUtImS14xqAEZeOAABAYgmoXq020.xlsx'
$ws.Range('C297').Value = 'http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=7ba8e15bffb24dffb06ea8ac7a8c0cb2&type=0'
$ws.Range('D297').Value = '中华人民共和国政府信息公开条例》'
$ws.Range('E297').Value = '《中华人民共和国政府信息公开条例》'
$ws.Range('F297').Value = '附件：巴中市平昌县五木镇政府信息主动公开基本目录.xlsx'
$ws.Range('G297').Value = 'http://www.scpc.gov.cn/ztzl/zfxxzdgkjbml/xzbsc/12617841.html'
$ws.Range('H297').Value = '巴中市平昌县五木镇政府信息主动公开基本目录'

# Row 298
$ws.Range('A298').Value = '2025-04-28 14:48'
$ws.Range('B298').Value = 'http://www.scpc.gov.cn/group3/M00/08/11/rBUtImIxibiANVSEAAE1Rnt3DRg20.docx'
$ws.Range('C298').Value = 'http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=0171c3b214e94f7d8118f44eb93c9d03&type=0'
$ws.Range('D298').Value = '社会和保障就业'
$ws.Range('E298').Value = '社会保障和就业'
$ws.Range('F298').Value = '平昌县林业局2019年度部门决算公开编制说明.docx.docx'
$ws.Range('G298').Value = 'http://www.scpc.gov.cn/public/6602341/12988791.html'
$ws.Range('H298').Value = '平昌县林业局2019年度部门决算公开'

# Row 299
$ws.Range('A299').Value = '2025-04-28 14:48'
$ws.Range('B299').Value = 'http://www.scpc.gov.cn/group3/M00/08/FD/rBUtImLCOE-ACmH0AADF69H66dc55.docx'
$ws.Range('C299').Value = 'http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=dc06e261b7f14ee0bf15f2175a30cb24&type=0'
$ws.Range('D299').Value = '抓住关键少数'
$ws.Range('E299').Value = '抓住“关键少数”'
$ws.Range('F299').Value = '平昌县纪委2018年部门决算公开编制说明.docx'
$ws.Range('G299').Value = 'http://www.scpc.gov.cn/public/6601841/12823521.html'
$ws.Range('H299').Value = '平昌县纪律检查委员会2018年度部门决算公开'

# Row 300
$ws.Range('A300').Value = '2025-04-28 14:49'
$ws.Range('B300').Value = 'http://www.scpc.gov.cn/group3/M00/07/A6/rBUtImHaoPqAeQCWAAI18HNAOO4495.pdf'
$ws.Range('C300').Value = 'http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=8503b427e2384f088efedc91da5a1ac7&type=0'
$ws.Range('D300').Value = '针炙'
$ws.Range('E300').Value = '针灸'
$ws.Range('F300').Value = '关于平昌县2021年下半年公开考试招聘卫生事业单位专业技术工作人员查分后公共科目笔试成绩（含政策性加分）岗位排名暨进入专业科目面试入围资格审查人员名单等相关事宜的公告'
$ws.Range('G300').Value = 'http://www.scpc.gov.cn/public/6602001/13804907.html'
$ws.Range('H300').Value = '关于平昌县2021年下半年公开考试招聘卫生事业单位专业技术工作人员查分后公共科目笔试成绩（含政策性加分...'

# Row 301
$ws.Range('A301').Value = '2025-04-28 14:49'
$ws.Range('B301').Value = 'http://www.scpc.gov.cn/group3/M00/08/17/rBUtImI4GRSAN1DrAAIEaATiT_8238.pdf'
$ws.Range('C301').Value = 'http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=905316a41d314c3591a91d22a1ad6d23&type=0'
$ws.Range('D301').Value = '针炙'
$ws.Range('E301').Value = '针灸'
$ws.Range('F301').Value = '关于平昌县2021年下半年公开招聘卫生事业单位专业技术工作人员拟聘人员的公示'
$ws.Range('G301').Value = 'http://www.scpc.gov.cn/zwgk/rsxx/13633451.html'
$ws.Range('H301').Value = '关于平昌县2021年下半年公开招聘卫生事业单位专业技术工作人员拟聘人员的公示'

# Row 302
$ws.Range('A302').Value = '2025-04-28 14:49'
$ws.Range('B302').Value = 'http://www.scpc.gov.cn/group3/M00/07/CB/rBUtImHyTRyAAl0hAAI0rhJ178A292.pdf'
$ws.Range('C302').Value = 'http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=2b4ffceb45f947c98a62972b3f6a7d66&type=0'
$ws.Range('D302').Value = '针炙'
$ws.Range('E302').Value = '针灸'
$ws.Range('F302').Value = '关于公布平昌县2021年下半年公开考试招聘卫生事业单位专业技术工作人员考试总成绩、岗位排名及进入体检人员名单的公告'
$ws.Range('G302').Value = 'http://www.scpc.gov.cn/public/6602001/13804904.html'
$ws.Range('H302').Value = '关于公布平昌县2021年下半年公开考试招聘卫生事业单位专业技术工作人员考试总成绩、岗位排名及进入体检人...'

# Row 303
$ws.Range('A303').Value = '2025-04-28 14:49'
$ws.Range('B303').Value = 'http://www.scpc.gov.cn/group3/M00/07/A6/rBUtImHaolOAFp50AAIzpsWy4js576.pdf'
$ws.Range('C303').Value = 'http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=675be85c769444bfb29b5aee711207f2&type=0'
$ws.Range('D303').Value = '针炙'
$ws.Range('E303').Value = '针灸'
$ws.Range('F303').Value = '关于平昌县2021年下半年公开考试招聘卫生事业单位专业技术工作人员查分后公共科目笔试成绩（含政策性加分）岗位排名暨进入专业科目面试入围资格审查人员名单等相关事宜的公告'
$ws.Range('G303').Value = 'http://www.scpc.gov.cn/public/6602001/13804907.html'
$ws.Range('H303').Value = '关于平昌县2021年下半年公开考试招聘卫生事业单位专业技术工作人员查分后公共科目笔试成绩（含政策性加分...'

# Row 304
$ws.Range('A304').Value = '2025-04-28 14:49'
$ws.Range('B304').Value = 'http://www.scpc.gov.cn/group3/M00/0E/E6/rBUtImUhAV-AWGGXAC8V976PGUI773.pdf'
$ws.Range('C304').Value = 'http://snapshot.ucap.com.cn/website/error-sensitive/index.html?id=0245404bc17441fd893e1f28b742e50d&type=0'
$ws.Range('D304').Value = '国士空间'
$ws.Range('E304').Value = '国土空间'
$ws.Range('F304').Value = '平安3井钻井工程--公示本.pdf'
$ws.Range('G304').Value = 'http://www.scpc.gov.cn/public/6602261/13893973.html'
$ws.Range('H304').Value = '关于2023年10月7日已受理建设项目环评文件公告(大庆油田有限责任公司)'
